$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the asset list (column A) and asset type (column B) while keeping
# quantities (column C) unchanged, reflecting the new portfolio for
# "Aula 2 Python Financas".

$ws.Range("A2").Value = "JHSF3"
$ws.Range("B2").Value = "Ação"

$ws.Range("A3").Value = "EGIE3"
$ws.Range("B3").Value = "Ação"

$ws.Range("A4").Value = "BBAS3"
$ws.Range("B4").Value = "Ação"

$ws.Range("A5").Value = "BBDC4"
$ws.Range("B5").Value = "Ação"

$ws.Range("A6").Value = "ITSA4"
$ws.Range("B6").Value = "Ação"

$ws.Range("A7").Value = "ABCB4"
$ws.Range("B7").Value = "Ação"

$ws.Range("A8").Value = "TRPL4"
$ws.Range("B8").Value = "Ação"

$ws.Range("A9").Value = "BBDC3"
$ws.Range("B9").Value = "Ação"

$ws.Range("A10").Value = "GNDI3"
$ws.Range("B10").Value = "Ação"

$ws.Range("A11").Value = "WEGE3"
$ws.Range("B11").Value = "Ação"

$ws.Range("A12").Value = "FLRY3"
$ws.Range("B12").Value = "Ação"

$ws.Range("A13").Value = "PASS3"
$ws.Range("B13").Value = "Ação"

$ws.Range("A14").Value = "BBSE3"
$ws.Range("B14").Value = "Ação"

$ws.Range("A15").Value = "RECR11"
$ws.Range("B15").Value = "FII"

$ws.Range("A10").Select()
